$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (shifts old N/O/P to O/P/Q), matching the
# width of the preceding column M so the new column's width serializes as "11".
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Activate this sheet and select J18, as in the edited workbook.
$ws.Activate()
$null = $ws.Range("J18").Select()
